$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab15")

# Update the source note: "April 2023" -> "October 2023" (IMF WEO database reference)
$ws.Range("A104").Value2 = "Source: Author's calculations based on IMF World Economic Outlook Database October 2023. Figures are adjusted so consumption, capital formation, changes in inventories and foreign balances as a percentage of GDP add up to 100."

# Refresh the aggregate (regional / income-group / total) rows with recalculated
# percent-of-GDP figures (columns C:H) following the IMF WEO October-2023 update.
    # Row 13
    $ws.Range("C13").Value2 = 78.505163009248307
    $ws.Range("D13").Value2 = 18.577486055557301
    $ws.Range("E13").Value2 = 0.77752802340920002
    $ws.Range("F13").Value2 = 2.1398229117852399
    $ws.Range("G13").Value2 = 33.3687189901663
    $ws.Range("H13").Value2 = 35.508541901951503
    # Row 23
    $ws.Range("C23").Value2 = 79.217475055368496
    $ws.Range("D23").Value2 = 17.664657419151599
    $ws.Range("E23").Value2 = 0.053038019608930002
    $ws.Range("F23").Value2 = 3.0648295058709598
    $ws.Range("G23").Value2 = 35.1284306370612
    $ws.Range("H23").Value2 = 38.193260142932203
    # Row 38
    $ws.Range("C38").Value2 = 82.608107524147002
    $ws.Range("D38").Value2 = 26.513217697470001
    $ws.Range("E38").Value2 = -0.33104655208669997
    $ws.Range("F38").Value2 = -8.7902786695303003
    $ws.Range("G38").Value2 = 24.3340102339887
    $ws.Range("H38").Value2 = 15.543731564458501
    # Row 45
    $ws.Range("C45").Value2 = 83.024019145858205
    $ws.Range("D45").Value2 = 18.921899364407999
    $ws.Range("E45").Value2 = 2.6377311668363199
    $ws.Range("F45").Value2 = -4.5836496771024002
    $ws.Range("G45").Value2 = 29.037986471862599
    $ws.Range("H45").Value2 = 24.454336794760099
    # Row 61
    $ws.Range("C61").Value2 = 81.3109691816165
    $ws.Range("D61").Value2 = 21.617331362391699
    $ws.Range("E61").Value2 = 0.84944230941719001
    $ws.Range("F61").Value2 = -3.7777428534254001
    $ws.Range("G61").Value2 = 23.391363470213101
    $ws.Range("H61").Value2 = 19.613620616787699
    # Row 62
    $ws.Range("C62").Value2 = 81.485089848854201
    $ws.Range("D62").Value2 = 20.612459349262
    $ws.Range("E62").Value2 = 1.2588632543260101
    $ws.Range("F62").Value2 = -3.3564124524421
    $ws.Range("G62").Value2 = 27.982187809780001
    $ws.Range("H62").Value2 = 24.625775357337901
    # Row 63
    $ws.Range("C63").Value2 = 70.780061594982598
    $ws.Range("D63").Value2 = 27.590699347626501
    $ws.Range("E63").Value2 = 1.30750430955534
    $ws.Range("F63").Value2 = 0.32173474783554001
    $ws.Range("G63").Value2 = 30.486891561539501
    $ws.Range("H63").Value2 = 30.759327244342799
    # Row 64
    $ws.Range("C64").Value2 = 81.331293087892206
    $ws.Range("D64").Value2 = 20.283935696399801
    $ws.Range("E64").Value2 = 0.31616236735950998
    $ws.Range("F64").Value2 = -1.9313911516515001
    $ws.Range("G64").Value2 = 30.7633814989786
    $ws.Range("H64").Value2 = 28.8319903473271
    # Row 65
    $ws.Range("C65").Value2 = 62.333681555912499
    $ws.Range("D65").Value2 = 36.318879054108002
    $ws.Range("E65").Value2 = 1.2589130121668
    $ws.Range("F65").Value2 = 0.088526377812699994
    $ws.Range("G65").Value2 = 23.6596735498598
    $ws.Range("H65").Value2 = 23.681389548938899
    # Row 66
    $ws.Range("C66").Value2 = 71.308590185659298
    $ws.Range("D66").Value2 = 27.246169724324901
    $ws.Range("E66").Value2 = 1.3051028036686001
    $ws.Range("F66").Value2 = 0.14013728634719
    $ws.Range("G66").Value2 = 30.362198570361102
    $ws.Range("H66").Value2 = 30.4539773877773
    # Row 67
    $ws.Range("C67").Value2 = 88.547748500693302
    $ws.Range("D67").Value2 = 16.6857375053157
    $ws.Range("E67").Value2 = 2.3810771209428299
    $ws.Range("F67").Value2 = -7.6145631269517997
    $ws.Range("G67").Value2 = 26.461137945188302
    $ws.Range("H67").Value2 = 18.846574818236501
    # Row 68
    $ws.Range("C68").Value2 = 85.229080731489304
    $ws.Range("D68").Value2 = 18.520924284648
    $ws.Range("E68").Value2 = 2.1025158013960801
    $ws.Range("F68").Value2 = -5.8525208175334003
    $ws.Range("G68").Value2 = 27.192974842802801
    $ws.Range("H68").Value2 = 21.340454025269398
    # Row 69
    $ws.Range("C69").Value2 = 82.531416543218995
    $ws.Range("D69").Value2 = 25.315218868201701
    $ws.Range("E69").Value2 = -0.4242350086144
    $ws.Range("F69").Value2 = -7.4224004028063
    $ws.Range("G69").Value2 = 26.325029044498802
    $ws.Range("H69").Value2 = 18.902628641692399
    # Row 70
    $ws.Range("C70").Value2 = 72.2728779727057
    $ws.Range("D70").Value2 = 20.233875275690401
    $ws.Range("E70").Value2 = 0.010898963646860001
    $ws.Range("F70").Value2 = 7.4823477879571101
    $ws.Range("G70").Value2 = 30.533453463132901
    $ws.Range("H70").Value2 = 38.015801251089997
    # Row 71
    $ws.Range("C71").Value2 = 81.3109691816165
    $ws.Range("D71").Value2 = 21.617331362391699
    $ws.Range("E71").Value2 = 0.84944230941719001
    $ws.Range("F71").Value2 = -3.7777428534254001
    $ws.Range("G71").Value2 = 23.391363470213101
    $ws.Range("H71").Value2 = 19.613620616787699
    # Row 72
    $ws.Range("C72").Value2 = 85.962199116529007
    $ws.Range("D72").Value2 = 23.648583498037599
    $ws.Range("E72").Value2 = 0.10227472358688
    $ws.Range("F72").Value2 = -9.7130573381533996
    $ws.Range("G72").Value2 = 21.759907236995801
    $ws.Range("H72").Value2 = 12.0468498988424
    # Row 73
    $ws.Range("C73").Value2 = 78.253528727374899
    $ws.Range("D73").Value2 = 20.7152542199366
    $ws.Range("E73").Value2 = 0.35941574161110001
    $ws.Range("F73").Value2 = 0.67180131107746999
    $ws.Range("G73").Value2 = 33.488235948013603
    $ws.Range("H73").Value2 = 34.160037259090998
    # Row 74
    $ws.Range("C74").Value2 = 73.012449867553201
    $ws.Range("D74").Value2 = 27.941988548845799
    $ws.Range("E74").Value2 = 0.33641836968078997
    $ws.Range("F74").Value2 = -1.2908567860797999
    $ws.Range("G74").Value2 = 39.605427993362298
    $ws.Range("H74").Value2 = 38.314571207282597
    # Row 75
    $ws.Range("C75").Value2 = 64.160365016026304
    $ws.Range("D75").Value2 = 29.897357863649798
    $ws.Range("E75").Value2 = 0.021786124243280001
    $ws.Range("F75").Value2 = 5.9204909960806198
    $ws.Range("G75").Value2 = 33.516792465614998
    $ws.Range("H75").Value2 = 39.4372834616956
    # Row 76
    $ws.Range("C76").Value2 = 68.103331877511394
    $ws.Range("D76").Value2 = 26.9959767662255
    $ws.Range("E76").Value2 = 2.0641224566127998
    $ws.Range("F76").Value2 = 2.8365688996503402
    $ws.Range("G76").Value2 = 50.435789190029297
    $ws.Range("H76").Value2 = 53.260400910542998
    # Row 77
    $ws.Range("C77").Value2 = 80.631057066369493
    $ws.Range("D77").Value2 = 19.681210705015399
    $ws.Range("E77").Value2 = 0.26863656005941999
    $ws.Range("F77").Value2 = -0.58090433144429998
    $ws.Range("G77").Value2 = 22.9814324458976
    $ws.Range("H77").Value2 = 22.400528114453301
    # Row 78
    $ws.Range("C78").Value2 = 74.258329749087096
    $ws.Range("D78").Value2 = 22.5806667216793
    $ws.Range("E78").Value2 = 1.97878214318177
    $ws.Range("F78").Value2 = 1.18222138605182
    $ws.Range("G78").Value2 = 54.549008974442899
    $ws.Range("H78").Value2 = 55.731230360494799
    # Row 79
    $ws.Range("C79").Value2 = 77.3155980445935
    $ws.Range("D79").Value2 = 22.715629722513
    $ws.Range("E79").Value2 = 1.32791258245413
    $ws.Range("F79").Value2 = -1.3591403495606
    $ws.Range("G79").Value2 = 34.7739631889701
    $ws.Range("H79").Value2 = 33.414823206107002
    # Row 80
    $ws.Range("C80").Value2 = 60.201084198860201
    $ws.Range("D80").Value2 = 29.549642323130801
    $ws.Range("E80").Value2 = -0.73547130189340004
    $ws.Range("F80").Value2 = 10.9847447799024
    $ws.Range("G80").Value2 = 26.706838020743401
    $ws.Range("H80").Value2 = 37.691582800645698
    # Row 81
    $ws.Range("C81").Value2 = 61.680880815556797
    $ws.Range("D81").Value2 = 22.332725236985102
    $ws.Range("E81").Value2 = 3.7998398685492498
    $ws.Range("F81").Value2 = 12.1865540789089
    $ws.Range("G81").Value2 = 25.346190360925501
    $ws.Range("H81").Value2 = 37.188568073675398
    # Row 82
    $ws.Range("C82").Value2 = 84.571946604215697
    $ws.Range("D82").Value2 = 19.3162839915564
    $ws.Range("E82").Value2 = 1.54810510092919
    $ws.Range("F82").Value2 = -5.4363356967012999
    $ws.Range("G82").Value2 = 28.167154031668701
    $ws.Range("H82").Value2 = 22.730818334967399
    # Row 83
    $ws.Range("C83").Value2 = 71.532876665813106
    $ws.Range("D83").Value2 = 28.025714570258099
    $ws.Range("E83").Value2 = 1.1013024749221001
    $ws.Range("F83").Value2 = -0.65989371099320004
    $ws.Range("G83").Value2 = 30.916246772062099
    $ws.Range("H83").Value2 = 30.222352239950599
    # Row 84
    $ws.Range("C84").Value2 = 85.177191756261394
    $ws.Range("D84").Value2 = 24.397045476809399
    $ws.Range("E84").Value2 = 0.28280812965315999
    $ws.Range("F84").Value2 = -9.8570453627239001
    $ws.Range("G84").Value2 = 33.0879080775183
    $ws.Range("H84").Value2 = 23.230862714794299
    # Row 86
    $ws.Range("C86").Value2 = 80.627756589462194
    $ws.Range("D86").Value2 = 20.897917241850202
    $ws.Range("E86").Value2 = 1.53127633348599
    $ws.Range("F86").Value2 = -3.0569501647983999
    $ws.Range("G86").Value2 = 25.638550963932801
    $ws.Range("H86").Value2 = 22.581600799134399
    # Row 87
    $ws.Range("C87").Value2 = 73.766747103404498
    $ws.Range("D87").Value2 = 28.732153053350899
    $ws.Range("E87").Value2 = 1.8955317918258801
    $ws.Range("F87").Value2 = -4.3944319485812002
    $ws.Range("G87").Value2 = 27.444438828795001
    $ws.Range("H87").Value2 = 22.642351129525199
    # Row 88
    $ws.Range("C88").Value2 = 81.840256886682496
    $ws.Range("D88").Value2 = 15.241155485426001
    $ws.Range("E88").Value2 = 0.94478644766040998
    $ws.Range("F88").Value2 = 1.9738011802311899
    $ws.Range("G88").Value2 = 33.855574077500798
    $ws.Range("H88").Value2 = 35.829375257732003
    # Row 89
    $ws.Range("C89").Value2 = 62.191613619463503
    $ws.Range("D89").Value2 = 34.061938358107099
    $ws.Range("E89").Value2 = 1.1375761791291901
    $ws.Range("F89").Value2 = 2.6088718433002298
    $ws.Range("G89").Value2 = 24.793238697664901
    $ws.Range("H89").Value2 = 27.327259136263098
    # Row 90
    $ws.Range("C90").Value2 = 75.892627901392601
    $ws.Range("D90").Value2 = 22.5679837489585
    $ws.Range("E90").Value2 = 1.24070343644199
    $ws.Range("F90").Value2 = 0.29868491320696999
    $ws.Range("G90").Value2 = 35.649294110735099
    $ws.Range("H90").Value2 = 35.947979400518498
    # Row 91
    $ws.Range("C91").Value2 = 78.675951683757901
    $ws.Range("D91").Value2 = 27.2920476078795
    $ws.Range("E91").Value2 = -0.18468255200479999
    $ws.Range("F91").Value2 = -5.7833167396325997
    $ws.Range("G91").Value2 = 31.267473436174701
    $ws.Range("H91").Value2 = 25.484156696542101
    # Row 92
    $ws.Range("C92").Value2 = 80.958375418329197
    $ws.Range("D92").Value2 = 30.692102936276299
    $ws.Range("E92").Value2 = 0.77887839135741999
    $ws.Range("F92").Value2 = -12.429356745963
    $ws.Range("G92").Value2 = 29.680895392622201
    $ws.Range("H92").Value2 = 17.251538646659299
    # Row 93
    $ws.Range("C93").Value2 = 89.561047414220596
    $ws.Range("D93").Value2 = 22.191375094863201
    $ws.Range("E93").Value2 = -0.053057416245400002
    $ws.Range("F93").Value2 = -11.699365092838001
    $ws.Range("G93").Value2 = 62.874953225289801
    $ws.Range("H93").Value2 = 51.175588132451601
    # Row 94
    $ws.Range("C94").Value2 = 56.321736486528899
    $ws.Range("D94").Value2 = 23.172650005412599
    $ws.Range("E94").Value2 = 0.99181163759275004
    $ws.Range("F94").Value2 = 19.513801870465802
    $ws.Range("G94").Value2 = 110.13184124147401
    $ws.Range("H94").Value2 = 129.64564311193999
    # Row 95
    $ws.Range("C95").Value2 = 83.620730474383393
    $ws.Range("D95").Value2 = 24.711997018853499
    $ws.Range("E95").Value2 = 0.38851897979536998
    $ws.Range("F95").Value2 = -8.7212464730323003
    $ws.Range("G95").Value2 = 28.606911399290301
    $ws.Range("H95").Value2 = 19.885664926258102
    # Row 96
    $ws.Range("C96").Value2 = 74.115697622998695
    $ws.Range("D96").Value2 = 24.417595728644301
    $ws.Range("E96").Value2 = 3.55715731203879
    $ws.Range("F96").Value2 = -2.0904506636817999
    $ws.Range("G96").Value2 = 37.897912056148797
    $ws.Range("H96").Value2 = 35.807461392467097
    # Row 97
    $ws.Range("C97").Value2 = 79.691397980803004
    $ws.Range("D97").Value2 = 24.5291070960319
    $ws.Range("E97").Value2 = -0.086635620693299997
    $ws.Range("F97").Value2 = -4.1338694561415998
    $ws.Range("G97").Value2 = 23.308418995661899
    $ws.Range("H97").Value2 = 19.174549539520299
    # Row 98
    $ws.Range("C98").Value2 = 78.604277470258694
    $ws.Range("D98").Value2 = 23.151717389687999
    $ws.Range("E98").Value2 = 4.6192246154434802
    $ws.Range("F98").Value2 = -6.3752194753901996
    $ws.Range("G98").Value2 = 23.914261042785999
    $ws.Range("H98").Value2 = 17.539041567395799
